$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Backlog - Open")

# Row 9 - SPIKE: Travel Planning Integration
$ws.Range("A9").Value = "SPIKE: Travel Planning Integration"
$ws.Range("B9").Value = "Icebox"
$ws.Range("C9").Value = 6
$ws.Range("D9").Value = 6
$ws.Range("E9").Value = "Low"
$ws.Range("F9").Value = 43477
$ws.Range("H9").Value = "User settles on a location.  Now how much would it cost to travel there?  We'll display cheapest 10%, median, and most expensive 10% of flights/hotels for that area"

# Row 10 - Filter undesirable locations
$ws.Range("A10").Value = "Filter undesirable locations"
$ws.Range("B10").Value = "Icebox"
$ws.Range("C10").Value = 7
$ws.Range("D10").Value = 3
$ws.Range("E10").Value = "Low"
$ws.Range("F10").Value = 43477
$ws.Range("H10").Value = "Does the user want to stay in the country?  Does the user want to go abroad?  Allow user to filter locations that are undesirable"

# Row 11 - SPIKE: Logo and User Interface
$ws.Range("A11").Value = "SPIKE: Logo and User Interface"
$ws.Range("B11").Value = "Icebox"
$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 3
$ws.Range("E11").Value = "Low"
$ws.Range("F11").Value = 43477
$ws.Range("H11").Value = "Once major functionality is implemented, make a nice, clear and usable interface"

# Update row 9 height to autofit the now-wrapped long note text
$ws.Rows("9").RowHeight = 31.5

# Update selection on Glossary sheet to A12 (without leaving it as the active tab)
$glossary = $wb.Worksheets.Item("Glossary")
$glossary.Range("A12").Select() | Out-Null

# Re-activate Backlog - Open and set its selection to A18 (final active sheet)
$ws.Activate()
$ws.Range("A18").Select() | Out-Null
